$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 185, pushing existing rows 185-304 down to 187-306.
$ws.Rows.Item(185).Resize(2).Insert()

# Copy the date cell style (numFmt) from the row right below (now row 187, originally row 185)
# down into the two newly inserted rows so the date values render/format correctly.
$ws.Range("D187").Copy()
$ws.Range("D185:D186").PasteSpecial(-4122)  # xlPasteFormats

# Populate new row 185
$ws.Cells.Item(185, 1).Value = 9
$ws.Cells.Item(185, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(185, 3).Value = "Metropolitana"
$ws.Cells.Item(185, 4).Value = 44438
$ws.Cells.Item(185, 5).Value = 13
$ws.Cells.Item(185, 6).Value = 100112040
$ws.Cells.Item(185, 7).Value = "Cilantro"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 43
$ws.Cells.Item(185, 11).Value = 10000
$ws.Cells.Item(185, 12).Value = 10000
$ws.Cells.Item(185, 13).Value = 10000
$ws.Cells.Item(185, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(185, 15).Value = "Región Metropolitana"
$ws.Cells.Item(185, 16).Value = 278
$ws.Cells.Item(185, 17).Value = 36
$ws.Cells.Item(185, 18).Value = "Hortaliza"

# Populate new row 186
$ws.Cells.Item(186, 1).Value = 9
$ws.Cells.Item(186, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(186, 3).Value = "Metropolitana"
$ws.Cells.Item(186, 4).Value = 44438
$ws.Cells.Item(186, 5).Value = 13
$ws.Cells.Item(186, 6).Value = 100112040
$ws.Cells.Item(186, 7).Value = "Cilantro"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 106
$ws.Cells.Item(186, 11).Value = 18000
$ws.Cells.Item(186, 12).Value = 20000
$ws.Cells.Item(186, 13).Value = 19000
$ws.Cells.Item(186, 14).Value = "`$/docena de atados"
$ws.Cells.Item(186, 15).Value = "Región Metropolitana"
$ws.Cells.Item(186, 16).Value = 6333
$ws.Cells.Item(186, 17).Value = 3
$ws.Cells.Item(186, 18).Value = "Hortaliza"
